# slides for April Teleconference
# Refresh the "datetimeFigureOut" footer date field (Insert > Header & Footer
# > Date and time) from 2/20/20 to 4/16/20 everywhere it is cached: the
# slide master, every slide layout, and the notes master.

$p = $ppt.ActivePresentation

$oldDate = "2/20/20"
$newDate = "4/16/20"

function Update-DateShape {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1 -and $sh.TextFrame.HasText -eq -1) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# 1. Slide master
Update-DateShape $p.SlideMaster.Shapes

# 2. Every slide layout that hangs off the (single) design/master
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Update-DateShape $layout.Shapes
}

# 3. Notes master
Update-DateShape $p.NotesMaster.Shapes
